$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Add explicit "none" borders (top/left/bottom/right) to the first table's
#    tblBorders, to match insideH/insideV which are already "none".
# ---------------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$borders = $t1.Borders
foreach ($idx in -1, -2, -3, -4) {
    $bd = $borders.Item($idx)
    # Order matters: reset width to 0 first, then set style to none, so the
    # resulting OOXML is <w:val="none" w:sz="0" .../> instead of sz="4".
    $bd.LineWidth = 0
    $bd.LineStyle = 0
}

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from right after "EXPERIENCE" to right
#    after "Competitive Programming, ".
# ---------------------------------------------------------------------------

# Remove the existing bookmark (currently sitting right after "EXPERIENCE").
if ($d.Bookmarks.Exists("_GoBack")) {
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()
}

# Locate the end of "Competitive Programming, ".
$findRange = $d.Content
$found = $findRange.Find.Execute("Competitive Programming, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPos = $findRange.End

# Adding a zero-length bookmark exactly at a paragraph's end-1 position
# (i.e. right before the paragraph mark) is flaky in this runtime, so we
# nudge the boundary out of the way: insert a throwaway character, add the
# bookmark next to it, then remove the throwaway character again. The
# bookmark stays anchored in the right spot.
$nudge = $d.Range($targetPos, $targetPos)
$nudge.InsertAfter("X")

$bmRange = $d.Range($targetPos, $targetPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($targetPos, $targetPos + 1).Delete()
